$d = $word.ActiveDocument

# Header date line
$d.Content.Find.Execute("2025-11-25 Tuesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-11-26 Wednesday", 2) | Out-Null

# Straightforward "old problem -> new problem" cell replacements.
$replacements = @(
    @("396×3=", "729×5="),
    @("750×3=", "868×8="),
    @("337×5=", "244×4="),
    @("686×9=", "332×2="),
    @("406×5=", "732×8="),
    @("584×5=", "738×6="),
    @("445×3=", "187×6="),
    @("748×4=", "867×4="),
    @("696×4=", "649×6="),
    @("255×9=", "543×7="),
    @("926×8=", "904×5="),
    @("496×8=", "441×3="),
    @("488×8=", "353×2="),
    @("252×9=", "346×9="),
    @("759×7=", "849×7="),
    @("406×2=", "544×5="),
    @("964×6=", "313×6="),
    @("374×2=", "863×4="),
    @("379×4=", "958×5="),
    @("325×9=", "563×6=")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# Last table row: problems were reshuffled (two new ones inserted, two old ones
# dropped) but the row still holds exactly five cells, so update the cell text
# in place rather than inserting/deleting cells.
$t = $d.Tables.Item(1)
$lastRow = $t.Rows.Count
$t.Cell($lastRow, 1).Range.Text = "671×5="
$t.Cell($lastRow, 2).Range.Text = "348×8="
$t.Cell($lastRow, 3).Range.Text = "719×8="
$t.Cell($lastRow, 4).Range.Text = "150×6="
$t.Cell($lastRow, 5).Range.Text = "414×9="
